$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AU1").Value = 0.99408914537612714
$ws.Range("BP1").Value = 0.79328199355768048
$ws.Range("O2").Value = 0.73609019227292061
$ws.Range("AL3").Value = 0.92144222316702473
$ws.Range("BI3").Value = 0.93930708634237203
$ws.Range("E4").Value = 0.80602507899354148
$ws.Range("F4").Value = 0.84312693604464151
$ws.Range("AL4").Value = 0.90270312090697802
$ws.Range("BH4").Value = 0.96388511247618402
$ws.Range("C5").Value = 0.72380411902231923
$ws.Range("AG5").Value = 0.57570338744458605
$ws.Range("W6").Value = 0.93998537394702297
$ws.Range("AW6").Value = 0.78897541064536991
$ws.Range("BG6").Value = 0.76675135725034038
$ws.Range("BN6").Value = 0.93481514647507002
$ws.Range("F8").Value = 0.98304685174672723
$ws.Range("AC8").Value = 0.87285642714262124
$ws.Range("AD9").Value = 0.94279007144098625
$ws.Range("I10").Value = 0.87749156407396178
$ws.Range("L10").Value = 0.62773752671535288
$ws.Range("M10").Value = 0.80766002652761237
$ws.Range("AE10").Value = 0.70934789974829582
$ws.Range("BJ10").Value = 0.84289913583344689
$ws.Range("AB11").Value = 0.65934579282230654
$ws.Range("AM11").Value = 0.96470932343928184
$ws.Range("BD11").Value = 0.59081382479117295
$ws.Range("BK11").Value = 0.78837470095441797
$ws.Range("AO12").Value = 0.84229249843077603
$ws.Range("L13").Value = 0.9229246775708102
$ws.Range("AJ13").Value = 0.94456882420726962
$ws.Range("BL14").Value = 0.79409885421296211
$ws.Range("AO15").Value = 0.91797754877577908
$ws.Range("AW17").Value = 0.67256115091972446
$ws.Range("P18").Value = 0.63008006429192176
$ws.Range("S18").Value = 0.93973055921909809
$ws.Range("AP19").Value = 0.80123415503464113
$ws.Range("BG19").Value = 0.60719814552017481
$ws.Range("Z20").Value = 0.85041281036098626
$ws.Range("BC20").Value = 0.91865599947410581
$ws.Range("BF20").Value = 0.97956543579313049
$ws.Range("W21").Value = 0.61143326650828245
$ws.Range("O22").Value = 0.91741521646419355
$ws.Range("Y22").Value = 0.99517084867771044
$ws.Range("AP22").Value = 0.90048413359830437
$ws.Range("R23").Value = 0.8679732945272749
$ws.Range("AY23").Value = 0.76266934265754338
$ws.Range("BK23").Value = 0.60975460616215604
$ws.Range("AN24").Value = 0.9874407272599558
$ws.Range("BA24").Value = 0.95243379670166139
$ws.Range("W25").Value = 0.88913810381887448
$ws.Range("AK26").Value = 0.7431541473424198
$ws.Range("AQ26").Value = 0.93197011388553619
$ws.Range("AY26").Value = 0.88647190012845756
$ws.Range("E27").Value = 0.92260976868060118
$ws.Range("U27").Value = 0.98996165036138362
$ws.Range("Y27").Value = 0.72101113306785081
$ws.Range("AC27").Value = 0.76258141617615594
$ws.Range("AR27").Value = 0.82324011508524619
$ws.Range("AL28").Value = 0.93096934061136116
$ws.Range("BH29").Value = 0.77093856763808744
$ws.Range("B30").Value = 0.86370830959680789
$ws.Range("S30").Value = 0.65553602971308611
$ws.Range("A31").Value = 0.75467294931717321
$ws.Range("B31").Value = 0.8807383871072566
$ws.Range("P31").Value = 0.6171099465825266
$ws.Range("AX31").Value = 0.99495172033211521
$ws.Range("K32").Value = 0.77276109181609209
$ws.Range("O32").Value = 0.84945897936879045
$ws.Range("P33").Value = 0.75078042590870198
$ws.Range("AQ33").Value = 0.83109640287862785
$ws.Range("AI34").Value = 0.95712290969037883
$ws.Range("AY35").Value = 0.92790726860407413
$ws.Range("BE35").Value = 0.91865708760199805
$ws.Range("I36").Value = 0.98253814556505537
$ws.Range("AH36").Value = 0.9070678836546151
$ws.Range("AL36").Value = 0.8836550893796391
$ws.Range("BH36").Value = 0.89485322768958397
$ws.Range("O37").Value = 0.89837969622824121
$ws.Range("AJ37").Value = 0.89200966163058271
$ws.Range("AU37").Value = 0.75150599212054381
$ws.Range("BB37").Value = 0.70457250428277307
$ws.Range("AP38").Value = 0.70601614507069876
$ws.Range("F39").Value = 0.89843871415690557
$ws.Range("AN39").Value = 0.99064001805651269
$ws.Range("AP39").Value = 0.8113006903838107
$ws.Range("AO40").Value = 0.91602589349449004
$ws.Range("AP40").Value = 0.85287234501690201
$ws.Range("BO40").Value = 0.66564654114361255
$ws.Range("G42").Value = 0.80483654249778858
$ws.Range("BC43").Value = 0.80800618491131138
$ws.Range("BG43").Value = 0.94440352905321467
$ws.Range("AG44").Value = 0.98336421075952796
$ws.Range("AT44").Value = 0.88356681031349282
$ws.Range("AJ45").Value = 0.94836626386905065
$ws.Range("Z46").Value = 0.60671863813065463
$ws.Range("AS46").Value = 0.67118819285672515
$ws.Range("Y47").Value = 0.92910118465759273
$ws.Range("BI47").Value = 0.88487312645051164
$ws.Range("AZ48").Value = 0.87599393831621319
$ws.Range("BO48").Value = 0.97377192648237898
$ws.Range("N49").Value = 0.87417671401344221
$ws.Range("AB49").Value = 0.91343480345629424
$ws.Range("AU49").Value = 0.95737047288307475
$ws.Range("AA50").Value = 0.64243117237778313
$ws.Range("AD52").Value = 0.91658016084258254
$ws.Range("AH52").Value = 0.87001477372812119
$ws.Range("AG53").Value = 0.91397885816504432
$ws.Range("AP54").Value = 0.66381974716012759
$ws.Range("BK54").Value = 0.83012429490563577
$ws.Range("L55").Value = 0.88502850376609699
$ws.Range("AF55").Value = 0.82805228729510749
$ws.Range("AT55").Value = 0.92410741726606749
$ws.Range("BO56").Value = 0.81336316011840826
$ws.Range("P57").Value = 0.97644466317135636
$ws.Range("AO58").Value = 0.75677982154795487
$ws.Range("N59").Value = 0.60382084454904761
$ws.Range("AS59").Value = 0.93142939006197012
$ws.Range("BE59").Value = 0.63681854767453361
$ws.Range("E60").Value = 0.64164498174105877
$ws.Range("BF60").Value = 0.92175888176378407
$ws.Range("S61").Value = 0.98866463690129225
$ws.Range("BA61").Value = 0.98855413290170202
$ws.Range("D62").Value = 0.96175945974449017
$ws.Range("AR62").Value = 0.77642600117939498
$ws.Range("T63").Value = 0.91074593094566803
$ws.Range("AU63").Value = 0.84386427813075204
$ws.Range("G65").Value = 0.85434866458719405
$ws.Range("AF65").Value = 0.94407773782295223
$ws.Range("BG65").Value = 0.7883431009860522
$ws.Range("Q66").Value = 0.99964256192350565
$ws.Range("AH66").Value = 0.78967377089661428
$ws.Range("BL66").Value = 0.84929801843254116
$ws.Range("AH67").Value = 0.81281876931468289
$ws.Range("H68").Value = 0.72011670836784591
$ws.Range("U68").Value = 0.80676393846636163
$ws.Range("Y68").Value = 0.78904429685137301
